$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel's
# smart "looks like a date" recognizer silently convert it into a date
# serial number (which would also stamp a date NumberFormat style onto
# the cell). We build the literal through a text formula, then collapse
# the formula down to its plain cached value via copy / paste-values.
function Set-TextValue {
    param($cell, [string]$text)
    $quote = [char]34
    $cell.Formula = "=" + $quote + $text + $quote
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# Row 2 (Coach 10805 / Axle 3 / Wheel 6): the flange-height alarm now
# predicts failure 5 days out, on 2017-04-12 instead of 2017-04-07.
$ws.Cells.Item(2, 7).Value = 5
Set-TextValue $ws.Cells.Item(2, 8) "2017-04-12"

# Row 3 becomes a brand-new record: a Gibson Ring inspection failure for
# Coach 10805 / Axle 2 / Wheel 4 (this pushes the former rows 3-5 down to
# 4-6, and the former row 6 drops off the bottom of the report).
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 4
$ws.Cells.Item(3, 4).Value = "75.5 in"
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 9).Value = "The Gibson Ring has Failed it's inspection"

# Row 4 (was Axle 4 / Wheel 8): shift the old row-3 data down + same
# flange-height alarm update as row 2.
$ws.Cells.Item(4, 3).Value = 8
$ws.Cells.Item(4, 5).Value = 30.1317
$ws.Cells.Item(4, 7).Value = 5
Set-TextValue $ws.Cells.Item(4, 8) "2017-04-12"

# Row 5 (was Axle 4 / Wheel 7): shift the old row-4 data down; predicted
# failure pushed out to 6 days / 2017-04-13.
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 7
$ws.Cells.Item(5, 5).Value = 30.064
$ws.Cells.Item(5, 7).Value = 6
Set-TextValue $ws.Cells.Item(5, 8) "2017-04-13"

# Row 6 (was Axle 3 / Wheel 5): shift the old row-5 data down; same as row 5.
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = 5
$ws.Cells.Item(6, 5).Value = 30.0098
$ws.Cells.Item(6, 7).Value = 6
Set-TextValue $ws.Cells.Item(6, 8) "2017-04-13"
